# The "Export" sheet re-sorted the balance list (descending by Saldo) and
# picked up a couple of new accounts. Net effect: only the first 28 data
# rows (rows 2-29, right under the "Conta/Nome/Saldo" header) actually
# change value/position; everything from row 30 down to the trailing
# "Filtros aplicados" footer is untouched, so we only touch rows 2-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Conta='004361159'; Nome='HFR';        Saldo=141845.61},
    @{Conta='004813166'; Nome='VENIA';      Saldo=99055.69},
    @{Conta='005690206'; Nome='KAUANNE';    Saldo=20020},
    @{Conta='004515548'; Nome='FLAVIA';     Saldo=9907.13},
    @{Conta='004886366'; Nome='RENATO';     Saldo=8806.33},
    @{Conta='004202332'; Nome='TATIANA';    Saldo=6604.48},
    @{Conta='004313254'; Nome='GUSTAVO';    Saldo=4292},
    @{Conta='004231371'; Nome='ADRIANO';    Saldo=3767.86},
    @{Conta='004368468'; Nome='AHMAD';      Saldo=2766.45},
    @{Conta='004213139'; Nome='LEONARDO';   Saldo=2609.8},
    @{Conta='004550415'; Nome='DIOGO';      Saldo=2155.74},
    @{Conta='004565146'; Nome='GUSTAVO';    Saldo=1980.94},
    @{Conta='004329030'; Nome='DANIELA';    Saldo=940.23},
    @{Conta='004392159'; Nome='RODRIGO';    Saldo=900.21},
    @{Conta='005696595'; Nome='CLUBE';      Saldo=752.05},
    @{Conta='001761119'; Nome='BLUEMETRIX'; Saldo=680.92},
    @{Conta='004855960'; Nome='CLERIA';     Saldo=556.35},
    @{Conta='004220849'; Nome='DULCE';      Saldo=503.59},
    @{Conta='008002502'; Nome='JORGEANA';   Saldo=500},
    @{Conta='005002457'; Nome='ROSANGELA';  Saldo=484.08},
    @{Conta='000806386'; Nome='FERNANDA';   Saldo=457.46},
    @{Conta='004432579'; Nome='ANA';        Saldo=446.18},
    @{Conta='004508516'; Nome='EDUARDO';    Saldo=364.49},
    @{Conta='004355790'; Nome='MINEIA';     Saldo=323.87},
    @{Conta='005040864'; Nome='ANDRE';      Saldo=279.96},
    @{Conta='004374891'; Nome='RODRIGO';    Saldo=273.15},
    @{Conta='004363260'; Nome='LARISSA';    Saldo=257.51},
    @{Conta='004515341'; Nome='BRUNO';      Saldo=235.12}
)

$r = 2
foreach ($row in $data) {
    # Prefix Conta with an apostrophe so Excel keeps the leading zeros as
    # text instead of coercing the digit string to a number.
    $ws.Cells.Item($r, 1).Value = "'" + $row.Conta
    $ws.Cells.Item($r, 2).Value = $row.Nome
    $ws.Cells.Item($r, 3).Value = $row.Saldo
    $r = $r + 1
}
